$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# NOTE: this runtime's PowerShell parser does not accept scientific-notation
# numeric literals (e.g. 1.23E-05), so every value below is written out in
# plain decimal form. The numeric value (and underlying IEEE-754 double) is
# identical either way.

# Row 2: D2/E2 6.881560286320602E-06 -> 5.405049002490183E-05
$ws.Range("D2").Value = 0.00005405049002490183
$ws.Range("E2").Value = 0.00005405049002490183

# Row 3: D3/E3 2.243249168317677E-55 -> 4.609702174208856E-77
$ws.Range("D3").Value = 0.00000000000000000000000000000000000000000000000000000000000000000000000000004609702174208856
$ws.Range("E3").Value = 0.00000000000000000000000000000000000000000000000000000000000000000000000000004609702174208856

# Row 4: D4/E4 3.759117941234606E-21 -> 1.227439719882987E-22
$ws.Range("D4").Value = 0.0000000000000000000001227439719882987
$ws.Range("E4").Value = 0.0000000000000000000001227439719882987

# Row 5: C5 Success TRUE -> FALSE ; D5/E5 1.479012157231247E-08 -> 0.9999985757577518
$ws.Range("C5").Value = $false
$ws.Range("D5").Value = 0.9999985757577518
$ws.Range("E5").Value = 0.9999985757577518

# Row 6: D6/E6 1 -> 0.9999999999991067
$ws.Range("D6").Value = 0.9999999999991067
$ws.Range("E6").Value = 0.9999999999991067

# Row 7: D7 0.9999999704596145 -> 0.9999999961869017 ; E7 2.954038547109405E-08 -> 3.813098281568728E-09
$ws.Range("D7").Value = 0.9999999961869017
$ws.Range("E7").Value = 0.000000003813098281568728

# Row 8: D8 0.9999999999550797 -> 0.9999999999999982 ; E8 4.492028971014861E-11 -> 1.77635683940025E-15
$ws.Range("D8").Value = 0.9999999999999982
$ws.Range("E8").Value = 0.00000000000000177635683940025

# Row 9: D9 0.9999999999812217 -> 0.9999999999690805 ; E9 1.87783122385099E-11 -> 3.091948919120568E-11
$ws.Range("D9").Value = 0.9999999999690805
$ws.Range("E9").Value = 0.00000000003091948919120568

# Row 10: D10 0.999984707895742 -> 0.9999969929415872 ; E10 1.529210425799121E-05 -> 3.007058412829444E-06
$ws.Range("D10").Value = 0.9999969929415872
$ws.Range("E10").Value = 0.000003007058412829444

# Row 11: D11 1 -> 0.9999999999999998 ; E11 0 -> 2.220446049250313E-16 ; F11 3.844633817672729 -> 4.12056827545166 ; G11 0.9 -> 0.8
$ws.Range("D11").Value = 0.9999999999999998
$ws.Range("E11").Value = 0.0000000000000002220446049250313
$ws.Range("F11").Value = 4.12056827545166
$ws.Range("G11").Value = 0.8
